$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Ph.D. Student, Cambridge University"
$ws.Range("A16").Value = "Ph.D. Plant Sciences "
$ws.Range("A17").Value = "M.Sc. Geography"
$ws.Range("A18").Value = "B.Sc. Environmental Science, Geographic Information Science (coop)"

$ws.Range("A18").Select()
